# Auto - Update data with bot!
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value  = "bbc-text.csv 데이터셋을 활용한 BBC 뉴스 아티클 카테고리 분류기 만들기"
$ws.Range("E4").Value  = "https://teddylee777.github.io/tensorflow/bbc-text-category-classification"

$ws.Range("D6").Value  = "[Python - 프로그래머스] 2023 KAKAO BLIND RECRUITMENT개인정보 수집 유효기간"
$ws.Range("E6").Value  = "https://leedakyeong.tistory.com/entry/Python-%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%98%EB%A8%B8%EC%8A%A4-2023-KAKAO-BLIND-RECRUITMENT%EA%B0%9C%EC%9D%B8%EC%A0%95%EB%B3%B4-%EC%88%98%EC%A7%91-%EC%9C%A0%ED%9A%A8%EA%B8%B0%EA%B0%84"

$ws.Range("D9").Value  = "공돌이들의 고질병 – “자동화”"
$ws.Range("E9").Value  = "https://blog.pabii.co.kr/engineer-sickness-automation/#utm_source=rss&utm_medium=rss&utm_campaign=engineer-sickness-automation"

$ws.Range("D12").Value = "“케라스 창시자에게 배우는 딥러닝 2판” 등 텐서플로 2.9.2 버전 업데이트 안내"
$ws.Range("E12").Value = "https://tensorflow.blog/2023/01/09/%ec%bc%80%eb%9d%bc%ec%8a%a4-%ec%b0%bd%ec%8b%9c%ec%9e%90%ec%97%90%ea%b2%8c-%eb%b0%b0%ec%9a%b0%eb%8a%94-%eb%94%a5%eb%9f%ac%eb%8b%9d-2%ed%8c%90-%eb%93%b1-%ed%85%90%ec%84%9c%ed%94%8c%eb%a1%9c-2-9-2/"

$ws.Range("D27").Value = "이루다 서버의 모니터링 스택을 소개합니다"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/spring-boot-monitoring-with-prometheus/"

$ws.Range("D28").Value = "[RVO] Reciprocal Velocity Obstacles"
$ws.Range("E28").Value = "https://ropiens.tistory.com/205"

$ws.Range("D32").Value = "Feature Interaction"
$ws.Range("E32").Value = "https://dodonam.tistory.com/399"

$ws.Range("D36").Value = "Clustering for Incomplete Time Series Data"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/392"

$ws.Range("D42").Value = "[임베디드-eclipse]Eclipse CDT 컴파일 옵션 설정"
$ws.Range("E42").Value = "https://kjk92.tistory.com/102"

$ws.Range("D50").Value = "가장 쉬운 Futsal Rules"
$ws.Range("E50").Value = "http://incredible.egloos.com/7576741"

$ws.Range("D51").Value = "[github actions] 깃허브 특정 브랜치에 push하는 순간 자동으로 도커 이미지 빌드해서 도커허브에 push하기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/github-actions-%EA%B9%83%ED%97%88%EB%B8%8C-%ED%8A%B9%EC%A0%95-%EB%B8%8C%EB%9E%9C%EC%B9%98%EC%97%90-push%ED%95%98%EB%8A%94-%EC%88%9C%EA%B0%84-%EC%9E%90%EB%8F%99%EC%9C%BC%EB%A1%9C-%EB%8F%84%EC%BB%A4-%EC%9D%B4%EB%AF%B8%EC%A7%80-%EB%B9%8C%EB%93%9C%ED%95%B4%EC%84%9C-%EB%8F%84%EC%BB%A4%ED%97%88%EB%B8%8C%EC%97%90-push%ED%95%98%EA%B8%B0"
